$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3528.25
$ws.Range("I6").Value = 292.375
$ws.Range("J6").Value = 10000
$ws.Range("K6").Value = 877.125
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = -765.125
$ws.Range("N6").Value = -30224
$ws.Range("H9").Value = 86
$ws.Range("I9").Value = 53.6
$ws.Range("J9").Value = 126.5
$ws.Range("K9").Value = 53.6
$ws.Range("L9").Value = 126.5
$ws.Range("M9").Value = 115.4
$ws.Range("N9").Value = -464.5
$ws.Range("H17").Value = 1198086.8
$ws.Range("J17").Value = 1227954.6
$ws.Range("L17").Value = 3683863.8
$ws.Range("N17").Value = -3684199.8
$ws.Range("H40").Value = 1699.8572
$ws.Range("I40").Value = 1599.8334
$ws.Range("J40").Value = 2300
$ws.Range("K40").Value = 1599.8334
$ws.Range("L40").Value = 2300
$ws.Range("M40").Value = -1424.8334
$ws.Range("N40").Value = -2650
$ws.Range("H51").Value = 7936.75
$ws.Range("J51").Value = 8713.571
$ws.Range("L51").Value = 8713.571
$ws.Range("N51").Value = -9681.571
$ws.Range("H74").Value = 7097.5
$ws.Range("I74").Value = 6130
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 6130
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = -5194
$ws.Range("N74").Value = -11872
$ws.Range("H77").Value = 7097.5
$ws.Range("I77").Value = 6130
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 30650
$ws.Range("L77").Value = 50000
$ws.Range("M77").Value = -25970
$ws.Range("N77").Value = -59360
$ws.Range("H87").Value = 55174.5
$ws.Range("J87").Value = 55174.5
$ws.Range("L87").Value = 55174.5
$ws.Range("N87").Value = -57670.5
$ws.Range("H90").Value = 55174.5
$ws.Range("J90").Value = 55174.5
$ws.Range("L90").Value = 165523.5
$ws.Range("N90").Value = -178003.5
$ws.Range("H98").Value = 1653.6786
$ws.Range("I98").Value = 505.15
$ws.Range("K98").Value = 505.15
$ws.Range("M98").Value = 992.85
$ws.Range("H122").Value = 1653.6786
$ws.Range("I122").Value = 505.15
$ws.Range("K122").Value = 1515.45
$ws.Range("M122").Value = 934.5500000000002
$ws.Range("H132").Value = 3199.5293
$ws.Range("I132").Value = 2579.6667
$ws.Range("J132").Value = 7848.5
$ws.Range("K132").Value = 7739.000100000001
$ws.Range("L132").Value = 23545.5
$ws.Range("M132").Value = -5209.000100000001
$ws.Range("N132").Value = -28605.5
$ws.Range("H137").Value = 7982.3257
$ws.Range("I137").Value = 4355.3335
$ws.Range("J137").Value = 12563.789
$ws.Range("K137").Value = 13066.0005
$ws.Range("L137").Value = 37691.367
$ws.Range("M137").Value = -10516.0005
$ws.Range("N137").Value = -42791.367
$ws.Range("H138").Value = 3136.5806
$ws.Range("I138").Value = 2815.3333
$ws.Range("J138").Value = 3811.2
$ws.Range("K138").Value = 8445.999899999999
$ws.Range("L138").Value = 11433.6
$ws.Range("M138").Value = -3305.999899999999
$ws.Range("N138").Value = -21713.6
$ws.Range("H140").Value = 499999
$ws.Range("J140").Value = 499999
$ws.Range("L140").Value = 499999
$ws.Range("N140").Value = -510359

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7288.524
$ws.Range("I45").Value = 11180.818
$ws.Range("K45").Value = 11180.818
$ws.Range("M45").Value = -10803.818
$ws.Range("H61").Value = 5177.914
$ws.Range("I61").Value = 3064.6072
$ws.Range("K61").Value = 3064.6072
$ws.Range("M61").Value = -2852.6072
$ws.Range("H114").Value = 70666.336
$ws.Range("J114").Value = 70666.336
$ws.Range("L114").Value = 70666.336
$ws.Range("N114").Value = -79344.336
$ws.Range("H124").Value = 29799.6
$ws.Range("J124").Value = 29799.6
$ws.Range("L124").Value = 29799.6
$ws.Range("N124").Value = -39619.6
$ws.Range("H136").Value = 5177.914
$ws.Range("I136").Value = 3064.6072
$ws.Range("K136").Value = 9193.821599999999
$ws.Range("M136").Value = -6643.821599999999
$ws.Range("H138").Value = 60000
$ws.Range("J138").Value = 60000
$ws.Range("L138").Value = 60000
$ws.Range("N138").Value = -70280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H88").Value = 17287
$ws.Range("J88").Value = 17287
$ws.Range("L88").Value = 17287
$ws.Range("N88").Value = -18099
$ws.Range("H91").Value = 17287
$ws.Range("J91").Value = 17287
$ws.Range("L91").Value = 17287
$ws.Range("N91").Value = -20095
$ws.Range("H99").Value = 2942.4614
$ws.Range("I99").Value = 1869.2858
$ws.Range("K99").Value = 1869.2858
$ws.Range("M99").Value = -371.2858000000001
$ws.Range("H105").Value = 2708.525
$ws.Range("I105").Value = 2969.2593
$ws.Range("K105").Value = 2969.2593
$ws.Range("M105").Value = -1222.2593
$ws.Range("H107").Value = 1550.9
$ws.Range("I107").Value = 1550.9
$ws.Range("K107").Value = 1550.9
$ws.Range("M107").Value = 369.0999999999999
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 7699.276
$ws.Range("I134").Value = 6585.5
$ws.Range("J134").Value = 11199.714
$ws.Range("K134").Value = 19756.5
$ws.Range("L134").Value = 33599.142
$ws.Range("M134").Value = -17221.5
$ws.Range("N134").Value = -38669.142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2050
$ws.Range("I3").Value = 2050
$ws.Range("K3").Value = 6150
$ws.Range("M3").Value = -6038
$ws.Range("H8").Value = 675.4286
$ws.Range("I8").Value = 675.4286
$ws.Range("K8").Value = 2026.2858
$ws.Range("M8").Value = -1887.2858
$ws.Range("H44").Value = 965
$ws.Range("I44").Value = 965
$ws.Range("K44").Value = 2895
$ws.Range("M44").Value = -2497
$ws.Range("H81").Value = 5881.75
$ws.Range("I81").Value = 6247.25
$ws.Range("J81").Value = 5516.25
$ws.Range("K81").Value = 18741.75
$ws.Range("L81").Value = 16548.75
$ws.Range("M81").Value = -17618.75
$ws.Range("N81").Value = -18794.75
$ws.Range("H84").Value = 5881.75
$ws.Range("I84").Value = 6247.25
$ws.Range("J84").Value = 5516.25
$ws.Range("K84").Value = 56225.25
$ws.Range("L84").Value = 49646.25
$ws.Range("M84").Value = -50609.25
$ws.Range("N84").Value = -60878.25
$ws.Range("H121").Value = 3601.818
$ws.Range("I121").Value = 3355.8572
$ws.Range("J121").Value = 4032.25
$ws.Range("K121").Value = 10067.5716
$ws.Range("L121").Value = 12096.75
$ws.Range("M121").Value = -8757.571599999999
$ws.Range("N121").Value = -14716.75
$ws.Range("H131").Value = 84985
$ws.Range("J131").Value = 3341.25
$ws.Range("L131").Value = 10023.75
$ws.Range("N131").Value = -20103.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 19068
$ws.Range("J15").Value = 19068
$ws.Range("L15").Value = 19068
$ws.Range("N15").Value = -19644
$ws.Range("H81").Value = 19068
$ws.Range("J81").Value = 19068
$ws.Range("L81").Value = 19068
$ws.Range("N81").Value = -21064
$ws.Range("H84").Value = 19068
$ws.Range("J84").Value = 19068
$ws.Range("L84").Value = 57204
$ws.Range("N84").Value = -67188
$ws.Range("H92").Value = 41998.332
$ws.Range("J92").Value = 41998.332
$ws.Range("L92").Value = 41998.332
$ws.Range("N92").Value = -45742.332
$ws.Range("H97").Value = 884.53845
$ws.Range("I97").Value = 943.0714
$ws.Range("J97").Value = 816.25
$ws.Range("K97").Value = 943.0714
$ws.Range("L97").Value = 816.25
$ws.Range("M97").Value = -447.0714
$ws.Range("N97").Value = -1808.25
$ws.Range("H98").Value = 32107.8
$ws.Range("J98").Value = 32107.8
$ws.Range("L98").Value = 32107.8
$ws.Range("N98").Value = -38097.8
$ws.Range("H113").Value = 162573.64
$ws.Range("I113").Value = 29783.111
$ws.Range("J113").Value = 401596.6
$ws.Range("K113").Value = 29783.111
$ws.Range("L113").Value = 401596.6
$ws.Range("M113").Value = -27613.111
$ws.Range("N113").Value = -405936.6
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 36472.875
$ws.Range("J74").Value = 36472.875
$ws.Range("L74").Value = 36472.875
$ws.Range("N74").Value = -38468.875
$ws.Range("H77").Value = 36472.875
$ws.Range("J77").Value = 36472.875
$ws.Range("L77").Value = 109418.625
$ws.Range("N77").Value = -119402.625
$ws.Range("H82").Value = 2228.65
$ws.Range("J82").Value = 2992.3333
$ws.Range("L82").Value = 2992.3333
$ws.Range("N82").Value = -3714.3333
$ws.Range("H85").Value = 2228.65
$ws.Range("J85").Value = 2992.3333
$ws.Range("L85").Value = 2992.3333
$ws.Range("N85").Value = -5488.3333
$ws.Range("H110").Value = 58429.668
$ws.Range("J110").Value = 58429.668
$ws.Range("L110").Value = 58429.668
$ws.Range("N110").Value = -66609.66800000001
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 44386
$ws.Range("J75").Value = 26629.5
$ws.Range("L75").Value = 26629.5
$ws.Range("N75").Value = -28501.5
$ws.Range("H78").Value = 44386
$ws.Range("J78").Value = 26629.5
$ws.Range("L78").Value = 79888.5
$ws.Range("N78").Value = -89248.5
